$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repulled data / mean calculation
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = -1
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 9
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = 3
$ws.Range("F24").Value = -3
$ws.Range("F25").Value = 9
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 11
$ws.Range("F30").Value = 3
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = -3
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = -2
$ws.Range("F36").Value = 2
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = -7
$ws.Range("F40").Value = 1
$ws.Range("F41").Value = -2
$ws.Range("F42").Value = -4
